$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1109.9
$ws.Range("I40").Value = 1112.375
$ws.Range("J40").Value = 1100
$ws.Range("K40").Value = 1112.375
$ws.Range("L40").Value = 1100
$ws.Range("M40").Value = -937.375
$ws.Range("N40").Value = -1450

$ws.Range("H112").Value = 1664.4445
$ws.Range("I112").Value = 750
$ws.Range("J112").Value = 2396
$ws.Range("K112").Value = 2250
$ws.Range("L112").Value = 7188
$ws.Range("M112").Value = -1142
$ws.Range("N112").Value = -9404

$ws.Range("H113").Value = 10001878
$ws.Range("I113").Value = 14287398
$ws.Range("J113").Value = 2333.3333
$ws.Range("K113").Value = 14287398
$ws.Range("L113").Value = 2333.3333
$ws.Range("M113").Value = -14284144
$ws.Range("N113").Value = -8841.3333

$ws.Range("H115").Value = 8653.210999999999
$ws.Range("I115").Value = 658.375
$ws.Range("J115").Value = 14467.637
$ws.Range("K115").Value = 1975.125
$ws.Range("L115").Value = 43402.911
$ws.Range("M115").Value = -408.125
$ws.Range("N115").Value = -46536.911

$ws.Range("H126").Value = 41134.285
$ws.Range("J126").Value = 41134.285
$ws.Range("L126").Value = 41134.285
$ws.Range("N126").Value = -51014.285

$ws.Range("H129").Value = 307508.25
$ws.Range("I129").Value = 59287.117
$ws.Range("J129").Value = 529600.8
$ws.Range("K129").Value = 177861.351
$ws.Range("L129").Value = 1588802.4
$ws.Range("M129").Value = -172861.351
$ws.Range("N129").Value = -1598802.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3012.32
$ws.Range("I32").Value = 2580.5454
$ws.Range("J32").Value = 6178.6665
$ws.Range("K32").Value = 2580.5454
$ws.Range("L32").Value = 6178.6665
$ws.Range("M32").Value = -2293.5454
$ws.Range("N32").Value = -6752.6665

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 31500
$ws.Range("J95").Value = 31500
$ws.Range("L95").Value = 31500
$ws.Range("N95").Value = -36992

$ws.Range("H134").Value = 12605632
$ws.Range("I134").Value = 12500569
$ws.Range("J134").Value = 14706882
$ws.Range("K134").Value = 37501707
$ws.Range("L134").Value = 44120646
$ws.Range("M134").Value = -37499172
$ws.Range("N134").Value = -44125716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10418208
$ws.Range("I31").Value = 22728312
$ws.Range("J31").Value = 1967.0769
$ws.Range("K31").Value = 22728312
$ws.Range("L31").Value = 1967.0769
$ws.Range("M31").Value = -22728017
$ws.Range("N31").Value = -2557.0769

$ws.Range("H34").Value = 10418208
$ws.Range("I34").Value = 22728312
$ws.Range("J34").Value = 1967.0769
$ws.Range("K34").Value = 22728312
$ws.Range("L34").Value = 1967.0769
$ws.Range("M34").Value = -22728110
$ws.Range("N34").Value = -2371.0769

$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H51").Value = 17049.25
$ws.Range("J51").Value = 17049.25
$ws.Range("L51").Value = 17049.25
$ws.Range("N51").Value = -18521.25

$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

$ws.Range("H59").Value = 30000
$ws.Range("J59").Value = 30000
$ws.Range("L59").Value = 30000
$ws.Range("N59").Value = -32290

$ws.Range("H61").Value = 17049.25
$ws.Range("J61").Value = 17049.25
$ws.Range("L61").Value = 17049.25
$ws.Range("N61").Value = -17745.25

$ws.Range("H62").Value = 3871.4285
$ws.Range("I62").Value = 2420
$ws.Range("J62").Value = 7500
$ws.Range("K62").Value = 2420
$ws.Range("L62").Value = 7500
$ws.Range("M62").Value = -1796
$ws.Range("N62").Value = -8748

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H65").Value = 3871.4285
$ws.Range("I65").Value = 2420
$ws.Range("J65").Value = 7500
$ws.Range("K65").Value = 12100
$ws.Range("L65").Value = 37500
$ws.Range("M65").Value = -8980
$ws.Range("N65").Value = -43740

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H69").Value = 18249.75
$ws.Range("I69").Value = 14333
$ws.Range("K69").Value = 14333
$ws.Range("M69").Value = -13584

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H72").Value = 18249.75
$ws.Range("I72").Value = 14333
$ws.Range("K72").Value = 42999
$ws.Range("M72").Value = -39255

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H75").Value = 40000
$ws.Range("J75").Value = 40000
$ws.Range("L75").Value = 40000
$ws.Range("N75").Value = -41996

$ws.Range("H78").Value = 40000
$ws.Range("J78").Value = 40000
$ws.Range("L78").Value = 120000
$ws.Range("N78").Value = -129984

$ws.Range("H129").Value = 49791.6
$ws.Range("J129").Value = 49791.6
$ws.Range("L129").Value = 49791.6
$ws.Range("N129").Value = -59791.6

$ws.Range("H131").Value = 15326
$ws.Range("J131").Value = 15326
$ws.Range("L131").Value = 15326
$ws.Range("N131").Value = -25406

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2564705.8
$ws.Range("I5").Value = 2564477.8
$ws.Range("J5").Value = 2564968.8
$ws.Range("K5").Value = 7693433.399999999
$ws.Range("L5").Value = 7694906.399999999
$ws.Range("M5").Value = -7693321.399999999
$ws.Range("N5").Value = -7695130.399999999

$ws.Range("H38").Value = 20041.6
$ws.Range("I38").Value = 33379.332
$ws.Range("J38").Value = 35
$ws.Range("K38").Value = 100137.996
$ws.Range("L38").Value = 105
$ws.Range("M38").Value = -99790.99600000001
$ws.Range("N38").Value = -799

$ws.Range("H107").Value = 475516.84
$ws.Range("I107").Value = 884385.0600000001
$ws.Range("J107").Value = 1229.72
$ws.Range("K107").Value = 2653155.18
$ws.Range("L107").Value = 3689.16
$ws.Range("M107").Value = -2651235.18
$ws.Range("N107").Value = -7529.16

$ws.Range("H131").Value = 48576
$ws.Range("I131").Value = 125492.5
$ws.Range("J131").Value = 1242.7693
$ws.Range("K131").Value = 376477.5
$ws.Range("L131").Value = 3728.3079
$ws.Range("M131").Value = -371437.5
$ws.Range("N131").Value = -13808.3079

$ws.Range("H135").Value = 2564705.8
$ws.Range("I135").Value = 2564477.8
$ws.Range("J135").Value = 2564968.8
$ws.Range("K135").Value = 23080300.2
$ws.Range("L135").Value = 23084719.2
$ws.Range("M135").Value = -23077765.2
$ws.Range("N135").Value = -23089789.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3948.125
$ws.Range("I102").Value = 4491.826
$ws.Range("J102").Value = 2558.6667
$ws.Range("K102").Value = 4491.826
$ws.Range("L102").Value = 2558.6667
$ws.Range("M102").Value = -2869.826
$ws.Range("N102").Value = -5802.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3475897.2
$ws.Range("I40").Value = 5558755.5
$ws.Range("J40").Value = 4466.9165
$ws.Range("K40").Value = 5558755.5
$ws.Range("L40").Value = 4466.9165
$ws.Range("M40").Value = -5558619.5
$ws.Range("N40").Value = -4738.9165

$ws.Range("H46").Value = 1118.88
$ws.Range("I46").Value = 1484.8235
$ws.Range("J46").Value = 341.25
$ws.Range("K46").Value = 1484.8235
$ws.Range("L46").Value = 341.25
$ws.Range("M46").Value = -1296.8235
$ws.Range("N46").Value = -717.25

$ws.Range("H82").Value = 4683.0586
$ws.Range("I82").Value = 1592.5555
$ws.Range("J82").Value = 8159.875
$ws.Range("K82").Value = 1592.5555
$ws.Range("L82").Value = 8159.875
$ws.Range("M82").Value = -1231.5555
$ws.Range("N82").Value = -8881.875

$ws.Range("H85").Value = 4683.0586
$ws.Range("I85").Value = 1592.5555
$ws.Range("J85").Value = 8159.875
$ws.Range("K85").Value = 1592.5555
$ws.Range("L85").Value = 8159.875
$ws.Range("M85").Value = -344.5554999999999
$ws.Range("N85").Value = -10655.875

$ws.Range("H136").Value = 1898211.8
$ws.Range("I136").Value = 1961452.8
$ws.Range("J136").Value = 980
$ws.Range("K136").Value = 5884358.4
$ws.Range("L136").Value = 2940
$ws.Range("M136").Value = -5881808.4
$ws.Range("N136").Value = -8040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2322.2222
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 2980
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 2980
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -5726

